$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: extend dimension by bulk-copying formats from the prior two 20-row blocks
$ws.Range("A1673:H1692").Copy()
$ws.Range("A1693:H1712").PasteSpecial(-4122)
$ws.Range("A1673:H1692").Copy()
$ws.Range("A1713:H1732").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 2: set cell values for the new rows
# Row 1693
$ws.Range("A1693").Value = 44009
$ws.Range("B1693").Value = "Hospital Universitario Miguel Servet"
$ws.Range("C1693").Value = 3
$ws.Range("D1693").Value = 2
$ws.Range("E1693").Value = "Zaragoza"
$ws.Range("F1693").Value = "Zaragoza"
$ws.Range("G1693").Value = 50297
$ws.Range("H1693").Value = "Fuente Aragón Hoy"

# Row 1694
$ws.Range("A1694").Value = 44009
$ws.Range("B1694").Value = "Hospital Clínico Universitario"
$ws.Range("C1694").Value = 30
$ws.Range("D1694").ClearContents()
$ws.Range("E1694").Value = "Zaragoza"
$ws.Range("F1694").Value = "Zaragoza"
$ws.Range("G1694").Value = 50297
$ws.Range("H1694").Value = "Fuente Aragón Hoy"

# Row 1695
$ws.Range("A1695").Value = 44009
$ws.Range("B1695").Value = "Hospital Royo Villanova"
$ws.Range("C1695").Value = 1
$ws.Range("D1695").ClearContents()
$ws.Range("E1695").Value = "Zaragoza"
$ws.Range("F1695").Value = "Zaragoza"
$ws.Range("G1695").Value = 50297
$ws.Range("H1695").Value = "Fuente Aragón Hoy"

# Row 1696
$ws.Range("A1696").Value = 44009
$ws.Range("B1696").Value = "Hospital Nuestra Señora de Gracia"
$ws.Range("C1696").Value = 1
$ws.Range("D1696").ClearContents()
$ws.Range("E1696").Value = "Zaragoza"
$ws.Range("F1696").Value = "Zaragoza"
$ws.Range("G1696").Value = 50297
$ws.Range("H1696").Value = "Fuente Aragón Hoy"

# Row 1697
$ws.Range("A1697").Value = 44009
$ws.Range("B1697").Value = "Hospital General de la Defensa"
$ws.Range("C1697").ClearContents()
$ws.Range("D1697").ClearContents()
$ws.Range("E1697").Value = "Zaragoza"
$ws.Range("F1697").Value = "Zaragoza"
$ws.Range("G1697").Value = 50297
$ws.Range("H1697").Value = "Fuente Aragón Hoy"

# Row 1698
$ws.Range("A1698").Value = 44009
$ws.Range("B1698").Value = "Hospital Obispo Polanco"
$ws.Range("C1698").Value = 1
$ws.Range("D1698").ClearContents()
$ws.Range("E1698").Value = "Teruel"
$ws.Range("F1698").Value = "Teruel"
$ws.Range("G1698").Value = 44216
$ws.Range("H1698").Value = "Fuente Aragón Hoy"

# Row 1699
$ws.Range("A1699").Value = 44009
$ws.Range("B1699").Value = "Hospital de Alcañiz"
$ws.Range("C1699").ClearContents()
$ws.Range("D1699").ClearContents()
$ws.Range("E1699").Value = "Alcañiz"
$ws.Range("F1699").Value = "Teruel"
$ws.Range("G1699").Value = 44013
$ws.Range("H1699").Value = "Fuente Aragón Hoy"

# Row 1700
$ws.Range("A1700").Value = 44009
$ws.Range("B1700").Value = "Hospital de Barbastro"
$ws.Range("C1700").Value = 9
$ws.Range("D1700").Value = 1
$ws.Range("E1700").Value = "Barbastro"
$ws.Range("F1700").Value = "Huesca"
$ws.Range("G1700").Value = 22048
$ws.Range("H1700").Value = "Fuente Aragón Hoy"

# Row 1701
$ws.Range("A1701").Value = 44009
$ws.Range("B1701").Value = "Hospital San Jorge"
$ws.Range("C1701").Value = 5
$ws.Range("D1701").ClearContents()
$ws.Range("E1701").Value = "Huesca"
$ws.Range("F1701").Value = "Huesca"
$ws.Range("G1701").Value = 22125
$ws.Range("H1701").Value = "Fuente Aragón Hoy"

# Row 1702
$ws.Range("A1702").Value = 44009
$ws.Range("B1702").Value = "Hospital Sagrado Corazón"
$ws.Range("C1702").ClearContents()
$ws.Range("D1702").ClearContents()
$ws.Range("E1702").Value = "Huesca"
$ws.Range("F1702").Value = "Huesca"
$ws.Range("G1702").Value = 22125
$ws.Range("H1702").Value = "Fuente Aragón Hoy"

# Row 1703
$ws.Range("A1703").Value = 44009
$ws.Range("B1703").Value = "Hospital Ernest Lluch"
$ws.Range("C1703").ClearContents()
$ws.Range("D1703").ClearContents()
$ws.Range("E1703").Value = "Calatayud"
$ws.Range("F1703").Value = "Zaragoza"
$ws.Range("G1703").Value = 50067
$ws.Range("H1703").Value = "Fuente Aragón Hoy"

# Row 1704
$ws.Range("A1704").Value = 44009
$ws.Range("B1704").Value = "Hospital San José"
$ws.Range("C1704").Value = 1
$ws.Range("D1704").ClearContents()
$ws.Range("E1704").Value = "Teruel"
$ws.Range("F1704").Value = "Teruel"
$ws.Range("G1704").Value = 44216
$ws.Range("H1704").Value = "Fuente Aragón Hoy"

# Row 1705
$ws.Range("A1705").Value = 44009
$ws.Range("B1705").Value = "Hospital Ejea – Cinco Villas"
$ws.Range("C1705").ClearContents()
$ws.Range("D1705").ClearContents()
$ws.Range("E1705").Value = "Ejea de los Caballeros"
$ws.Range("F1705").Value = "Zaragoza"
$ws.Range("G1705").Value = 50095
$ws.Range("H1705").Value = "Fuente Aragón Hoy"

# Row 1706
$ws.Range("A1706").Value = 44009
$ws.Range("B1706").Value = "MAZ"
$ws.Range("C1706").ClearContents()
$ws.Range("D1706").ClearContents()
$ws.Range("E1706").Value = "Zaragoza"
$ws.Range("F1706").Value = "Zaragoza"
$ws.Range("G1706").Value = 50297
$ws.Range("H1706").Value = "Fuente Aragón Hoy"

# Row 1707
$ws.Range("A1707").Value = 44009
$ws.Range("B1707").Value = "Hospital Viamed Montecanal"
$ws.Range("C1707").ClearContents()
$ws.Range("D1707").ClearContents()
$ws.Range("E1707").Value = "Zaragoza"
$ws.Range("F1707").Value = "Zaragoza"
$ws.Range("G1707").Value = 50297
$ws.Range("H1707").Value = "Fuente Aragón Hoy"

# Row 1708
$ws.Range("A1708").Value = 44009
$ws.Range("B1708").Value = "Clínica Montpellier"
$ws.Range("C1708").ClearContents()
$ws.Range("D1708").ClearContents()
$ws.Range("E1708").Value = "Zaragoza"
$ws.Range("F1708").Value = "Zaragoza"
$ws.Range("G1708").Value = 50297
$ws.Range("H1708").Value = "Fuente Aragón Hoy"

# Row 1709
$ws.Range("A1709").Value = 44009
$ws.Range("B1709").Value = "Hospital Quirón"
$ws.Range("C1709").ClearContents()
$ws.Range("D1709").ClearContents()
$ws.Range("E1709").Value = "Zaragoza"
$ws.Range("F1709").Value = "Zaragoza"
$ws.Range("G1709").Value = 50297
$ws.Range("H1709").Value = "Fuente Aragón Hoy"

# Row 1710
$ws.Range("A1710").Value = 44009
$ws.Range("B1710").Value = "Hospital San Juan de Dios de Zaragoza"
$ws.Range("C1710").ClearContents()
$ws.Range("D1710").ClearContents()
$ws.Range("E1710").Value = "Zaragoza"
$ws.Range("F1710").Value = "Zaragoza"
$ws.Range("G1710").Value = 50297
$ws.Range("H1710").Value = "Fuente Aragón Hoy"

# Row 1711
$ws.Range("A1711").Value = 44009
$ws.Range("B1711").Value = "Clínica Viamed Santiago"
$ws.Range("C1711").ClearContents()
$ws.Range("D1711").ClearContents()
$ws.Range("E1711").Value = "Huesca"
$ws.Range("F1711").Value = "Huesca"
$ws.Range("G1711").Value = 22125
$ws.Range("H1711").Value = "Fuente Aragón Hoy"

# Row 1712
$ws.Range("A1712").Value = 44009
$ws.Range("B1712").Value = "Clínica El Pilar"
$ws.Range("C1712").ClearContents()
$ws.Range("D1712").ClearContents()
$ws.Range("E1712").Value = "Zaragoza"
$ws.Range("F1712").Value = "Zaragoza"
$ws.Range("G1712").Value = 50297
$ws.Range("H1712").Value = "Fuente Aragón Hoy"

# Row 1713
$ws.Range("A1713").Value = 44010
$ws.Range("B1713").Value = "Hospital Universitario Miguel Servet"
$ws.Range("C1713").Value = 3
$ws.Range("D1713").Value = 2
$ws.Range("E1713").Value = "Zaragoza"
$ws.Range("F1713").Value = "Zaragoza"
$ws.Range("G1713").Value = 50297
$ws.Range("H1713").Value = "Fuente Aragón Hoy"

# Row 1714
$ws.Range("A1714").Value = 44010
$ws.Range("B1714").Value = "Hospital Clínico Universitario"
$ws.Range("C1714").Value = 22
$ws.Range("D1714").ClearContents()
$ws.Range("E1714").Value = "Zaragoza"
$ws.Range("F1714").Value = "Zaragoza"
$ws.Range("G1714").Value = 50297
$ws.Range("H1714").Value = "Fuente Aragón Hoy"

# Row 1715
$ws.Range("A1715").Value = 44010
$ws.Range("B1715").Value = "Hospital Royo Villanova"
$ws.Range("C1715").ClearContents()
$ws.Range("D1715").ClearContents()
$ws.Range("E1715").Value = "Zaragoza"
$ws.Range("F1715").Value = "Zaragoza"
$ws.Range("G1715").Value = 50297
$ws.Range("H1715").Value = "Fuente Aragón Hoy"

# Row 1716
$ws.Range("A1716").Value = 44010
$ws.Range("B1716").Value = "Hospital Nuestra Señora de Gracia"
$ws.Range("C1716").Value = 1
$ws.Range("D1716").ClearContents()
$ws.Range("E1716").Value = "Zaragoza"
$ws.Range("F1716").Value = "Zaragoza"
$ws.Range("G1716").Value = 50297
$ws.Range("H1716").Value = "Fuente Aragón Hoy"

# Row 1717
$ws.Range("A1717").Value = 44010
$ws.Range("B1717").Value = "Hospital General de la Defensa"
$ws.Range("C1717").ClearContents()
$ws.Range("D1717").ClearContents()
$ws.Range("E1717").Value = "Zaragoza"
$ws.Range("F1717").Value = "Zaragoza"
$ws.Range("G1717").Value = 50297
$ws.Range("H1717").Value = "Fuente Aragón Hoy"

# Row 1718
$ws.Range("A1718").Value = 44010
$ws.Range("B1718").Value = "Hospital Obispo Polanco"
$ws.Range("C1718").Value = 1
$ws.Range("D1718").ClearContents()
$ws.Range("E1718").Value = "Teruel"
$ws.Range("F1718").Value = "Teruel"
$ws.Range("G1718").Value = 44216
$ws.Range("H1718").Value = "Fuente Aragón Hoy"

# Row 1719
$ws.Range("A1719").Value = 44010
$ws.Range("B1719").Value = "Hospital de Alcañiz"
$ws.Range("C1719").ClearContents()
$ws.Range("D1719").ClearContents()
$ws.Range("E1719").Value = "Alcañiz"
$ws.Range("F1719").Value = "Teruel"
$ws.Range("G1719").Value = 44013
$ws.Range("H1719").Value = "Fuente Aragón Hoy"

# Row 1720
$ws.Range("A1720").Value = 44010
$ws.Range("B1720").Value = "Hospital de Barbastro"
$ws.Range("C1720").Value = 7
$ws.Range("D1720").Value = 1
$ws.Range("E1720").Value = "Barbastro"
$ws.Range("F1720").Value = "Huesca"
$ws.Range("G1720").Value = 22048
$ws.Range("H1720").Value = "Fuente Aragón Hoy"

# Row 1721
$ws.Range("A1721").Value = 44010
$ws.Range("B1721").Value = "Hospital San Jorge"
$ws.Range("C1721").Value = 4
$ws.Range("D1721").ClearContents()
$ws.Range("E1721").Value = "Huesca"
$ws.Range("F1721").Value = "Huesca"
$ws.Range("G1721").Value = 22125
$ws.Range("H1721").Value = "Fuente Aragón Hoy"

# Row 1722
$ws.Range("A1722").Value = 44010
$ws.Range("B1722").Value = "Hospital Sagrado Corazón"
$ws.Range("C1722").ClearContents()
$ws.Range("D1722").ClearContents()
$ws.Range("E1722").Value = "Huesca"
$ws.Range("F1722").Value = "Huesca"
$ws.Range("G1722").Value = 22125
$ws.Range("H1722").Value = "Fuente Aragón Hoy"

# Row 1723
$ws.Range("A1723").Value = 44010
$ws.Range("B1723").Value = "Hospital Ernest Lluch"
$ws.Range("C1723").ClearContents()
$ws.Range("D1723").ClearContents()
$ws.Range("E1723").Value = "Calatayud"
$ws.Range("F1723").Value = "Zaragoza"
$ws.Range("G1723").Value = 50067
$ws.Range("H1723").Value = "Fuente Aragón Hoy"

# Row 1724
$ws.Range("A1724").Value = 44010
$ws.Range("B1724").Value = "Hospital San José"
$ws.Range("C1724").Value = 1
$ws.Range("D1724").ClearContents()
$ws.Range("E1724").Value = "Teruel"
$ws.Range("F1724").Value = "Teruel"
$ws.Range("G1724").Value = 44216
$ws.Range("H1724").Value = "Fuente Aragón Hoy"

# Row 1725
$ws.Range("A1725").Value = 44010
$ws.Range("B1725").Value = "Hospital Ejea – Cinco Villas"
$ws.Range("C1725").ClearContents()
$ws.Range("D1725").ClearContents()
$ws.Range("E1725").Value = "Ejea de los Caballeros"
$ws.Range("F1725").Value = "Zaragoza"
$ws.Range("G1725").Value = 50095
$ws.Range("H1725").Value = "Fuente Aragón Hoy"

# Row 1726
$ws.Range("A1726").Value = 44010
$ws.Range("B1726").Value = "MAZ"
$ws.Range("C1726").ClearContents()
$ws.Range("D1726").ClearContents()
$ws.Range("E1726").Value = "Zaragoza"
$ws.Range("F1726").Value = "Zaragoza"
$ws.Range("G1726").Value = 50297
$ws.Range("H1726").Value = "Fuente Aragón Hoy"

# Row 1727
$ws.Range("A1727").Value = 44010
$ws.Range("B1727").Value = "Hospital Viamed Montecanal"
$ws.Range("C1727").ClearContents()
$ws.Range("D1727").ClearContents()
$ws.Range("E1727").Value = "Zaragoza"
$ws.Range("F1727").Value = "Zaragoza"
$ws.Range("G1727").Value = 50297
$ws.Range("H1727").Value = "Fuente Aragón Hoy"

# Row 1728
$ws.Range("A1728").Value = 44010
$ws.Range("B1728").Value = "Clínica Montpellier"
$ws.Range("C1728").Value = 3
$ws.Range("D1728").ClearContents()
$ws.Range("E1728").Value = "Zaragoza"
$ws.Range("F1728").Value = "Zaragoza"
$ws.Range("G1728").Value = 50297
$ws.Range("H1728").Value = "Fuente Aragón Hoy"

# Row 1729
$ws.Range("A1729").Value = 44010
$ws.Range("B1729").Value = "Hospital Quirón"
$ws.Range("C1729").ClearContents()
$ws.Range("D1729").ClearContents()
$ws.Range("E1729").Value = "Zaragoza"
$ws.Range("F1729").Value = "Zaragoza"
$ws.Range("G1729").Value = 50297
$ws.Range("H1729").Value = "Fuente Aragón Hoy"

# Row 1730
$ws.Range("A1730").Value = 44010
$ws.Range("B1730").Value = "Hospital San Juan de Dios de Zaragoza"
$ws.Range("C1730").ClearContents()
$ws.Range("D1730").ClearContents()
$ws.Range("E1730").Value = "Zaragoza"
$ws.Range("F1730").Value = "Zaragoza"
$ws.Range("G1730").Value = 50297
$ws.Range("H1730").Value = "Fuente Aragón Hoy"

# Row 1731
$ws.Range("A1731").Value = 44010
$ws.Range("B1731").Value = "Clínica Viamed Santiago"
$ws.Range("C1731").ClearContents()
$ws.Range("D1731").ClearContents()
$ws.Range("E1731").Value = "Huesca"
$ws.Range("F1731").Value = "Huesca"
$ws.Range("G1731").Value = 22125
$ws.Range("H1731").Value = "Fuente Aragón Hoy"

# Row 1732
$ws.Range("A1732").Value = 44010
$ws.Range("B1732").Value = "Clínica El Pilar"
$ws.Range("C1732").ClearContents()
$ws.Range("D1732").ClearContents()
$ws.Range("E1732").Value = "Zaragoza"
$ws.Range("F1732").Value = "Zaragoza"
$ws.Range("G1732").Value = 50297
$ws.Range("H1732").Value = "Fuente Aragón Hoy"
